{"js": "const replacements = [\n  [\n    \"Ativa\u00e7\u00e3o: 01/01/2022\",\n    \"Ativa\u00e7\u00e3o: 01/01/2024\"\n  ],\n  [\n    \"Curso (semestre ideal): EQD (9), EQN (10)\",\n    \"Curso (semestre ideal): EQN (10)\"\n  ],\n  [\n    \"Conferir aos alunos uma no\u00e7\u00e3o pr\u00e1tica das aplica\u00e7\u00f5es, \u00e0 escala industrial, de processos e produtos estudados em disciplinas de qu\u00edmica org\u00e2nica e ainda uma vis\u00e3o global das mat\u00e9rias primas mais importantes na \u00e1rea da ind\u00fastria qu\u00edmica de base org\u00e2nica.\",\n    \"Conferir aos alunos aplica\u00e7\u00f5es em escala industrial, de processos e produtos estudados em disciplinas de qu\u00edmica org\u00e2nica e uma vis\u00e3o global das mat\u00e9rias primas mais importantes na \u00e1rea da ind\u00fastria qu\u00edmica de base org\u00e2nica.\"\n  ],\n  [\n    \"Check the students a practical notion of applications on an industrial scale, processes and products studied in organic chemistry disciplines and also an overview of the most important raw materials in the chemical industry of organic base\",\n    \"Check the students students with applications on an industrial scale of processes and products studied in organic chemistry disciplines and an overview of the most important raw materials in the chemical industry of organic base.\"\n  ],\n  [\n    \"Petr\u00f3leo, G\u00e1s Natural e Petroqu\u00edmica; Qu\u00edmica Fina; Processos Unit\u00e1rios Org\u00e2nicos: Nitra\u00e7\u00e3o; Esterifica\u00e7\u00e3o; Alquila\u00e7\u00e3o e Acila\u00e7\u00e3o; Hidrogena\u00e7\u00e3o; Sulfona\u00e7\u00e3o/Sulfata\u00e7\u00e3o; Oxida\u00e7\u00e3o.\",\n    \"Qu\u00edmica Fina,Petr\u00f3leo, G\u00e1s Natural e Petroqu\u00edmica; Processos Unit\u00e1rios Org\u00e2nicos:Alquila\u00e7\u00e3o e Acila\u00e7\u00e3o; Hidrogena\u00e7\u00e3o e Desidrogena\u00e7\u00e3o,Halogena\u00e7\u00e3o, Esterifica\u00e7\u00e3o; Sulfona\u00e7\u00e3o/Sulfata\u00e7\u00e3o; Oxida\u00e7\u00e3o.\"\n  ],\n  [\n    \"Oil, Natural Gas and Petrochemicals; Fine Chemistry; Organic Unit Processes: Nitration; Esterification; Alkylation and Acylation; Hydrogenation; Sulphonation/Sulfation; Oxidation.\",\n    \"Fine Chemistry,Oil, Natural Gas and Petrochemicals; Organic Unit Processes: Alkylation and Acylation; Hydrogenation,Halogenation, Esterification; Sulphonation/Sulfation; Oxidation.\"\n  ],\n  [\n    \"Petr\u00f3leo, G\u00e1s Natural e Petroqu\u00edmica\u037e 2- Qu\u00edmica Fina: Caracter\u00edsticas, Qu\u00edmica Fina X Qu\u00edmica de Base, Principais Segmentos (Defensivos Agr\u00edcolas, F\u00e1rmacos, Catalisadores, Corantes e Pigmentos, Especialidades); 3- Processos Unit\u00e1rios Org\u00e2nicos: 3.1- Nitra\u00e7\u00e3o; 3.2- Esterifica\u00e7\u00e3o; 3.3- Alquila\u00e7\u00e3o e Acila\u00e7\u00e3o; 3.4- Hidrogena\u00e7\u00e3o; 3.5- Sulfona\u00e7\u00e3o/Sulfata\u00e7\u00e3o; 3.6- Oxida\u00e7\u00e3o.\",\n    \"1-Qu\u00edmica Fina: 1.1- Caracter\u00edsticas, 1.2- Principais Segmentos (Defensivos Agr\u00edcolas, F\u00e1rmacos, Catalisadores, Corantes e Pigmentos, Especialidades), 1.3- Qu\u00edmica Fina X Qu\u00edmica de Base, 2- Petr\u00f3leo, G\u00e1s Natural e Petroqu\u00edmica; 3- Processos Unit\u00e1rios Org\u00e2nicos: 3.1- Alquila\u00e7\u00e3o e Acila\u00e7\u00e3o; 3.2- Hidrogena\u00e7\u00e3o e Desidogena\u00e7\u00e3o;3.2.1 Processos Oxo, 3.2.2- Amino;3.3- Halogena\u00e7\u00e3o; 3.4- Esterifica\u00e7\u00e3o; 3.5- Sulfona\u00e7\u00e3o/Sulfata\u00e7\u00e3o; 3.6- Oxida\u00e7\u00e3o.\"\n  ],\n  [\n    \"1- Oil, Natural Gas and Petrochemicals 2- Fine Chemistry: Characteristics, Fine Chemistry vs. Basic Chemistry, Main Segments (Pesticides, Drugs, Catalysts, Dyes and Pigments, Specialties); 3- Organic Unit Processes: 3.1- Nitration; 3.2- Esterification; 3.3- Alkylation and Acylation; 3.4- Hydrogenation; 3.5- Sulphonation/Sulfation; 3.6- Oxidation.\",\n    \"1- Fine Chemicals: 1.1- Characteristics, 1.2- Main Segments (Agricultural Defensives, Pharmaceuticals, Catalysts, Dyes and Pigments, Specialties), 1.3- Fine Chemicals X Basic Chemicals, 2- Oil, Natural Gas and Petrochemicals; 3- Unitary Processes Organic: 3.1- Alkylation and Acylation; 3.2- Hydrogenation and Dehydrogenation; 3.2.1 Oxo processes, 3.2.2- Amino; 3.3- Halogenation; 3.4- Esterification; 3.5- Sulfonation/Sulfation; 3.6- Oxidation.\"\n  ],\n  [\n    \"Provas em sala, entrega de exerc\u00edcios ou casos pr\u00e1ticos elaborados fora de sala de aula.\",\n    \"A nota (NOTA) ser\u00e1 composta por uma destas op\u00e7\u00f5es: prova em sala, apresenta\u00e7\u00f5es em sala, entrega de exerc\u00edcios ou casos pr\u00e1ticos elaborados fora de sala de aula. A estas op\u00e7\u00f5es ser\u00e1 incorporado,para cada aluno, seu respectivo percentual defrequ\u00eancia no c\u00e1lculo da nota final (NF), conforme a f\u00f3rmula explicitada abaixo:NF = NOTA x %FREQ.\"\n  ],\n  [\n    \"Frequ\u00eancia m\u00ednima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recupera\u00e7\u00e3o.\",\n    \"Frequ\u00eancia m\u00ednima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita aplica\u00e7\u00e3o de prova escrita de recupera\u00e7\u00e3o valendo 10,00 pontos.\"\n  ],\n  [\n    \"Ullmann\u2019s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim \u037e New York : WileyVCH, 2011.Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.Manual econ\u00f4mico da ind\u00fastria qu\u00edmica - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Cama\u00e7ari: CEPED, 2007.Shreve, R. Norris\u037e BRINK JR., J. A. Ind\u00fastrias de processos qu\u00edmicos. Tradu\u00e7\u00e3o de Hor\u00e1cio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.T.W. Graham Solomons, Craig B. Fryhle Hoboken, NJ. Organic chemistry; John Wiley, 9th ed; c2008.Revistas:Brazilian Journal of Chemical Engineering, S\u00e3o Paulo, SP: Brazilian Society of Chemical Engineering, v. 11, n. 1, 1995-;Qu\u00edmica & Derivados, S\u00e3o Paulo, SP: QD, v. 1, n. 1, nov. 1965-; Dispon\u00edvel em: http://www.quimica.com.br/pquimica/category/revista/BiodieselBR. Curitiba: BiodieselBR, v.1, n.1, out/dez.2007 -;Petr\u00f3leo & Energia, S\u00e3o Paulo, SP: , v. 1, n. 1, ; Dispon\u00edvel em: http://www.petroleoenergia.com.br/petroleo/category/revista-petroleo-e-energia/.\",\n    \"Livros:Ullmann\u2019s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim ; New York : WileyVCH, 2011.Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.Manual Econ\u00f4mico da Ind\u00fastria Qu\u00edmica - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Cama\u00e7ari: CEPED, 2007.Shreve, R. Norris; BRINK JR., J. A. Ind\u00fastrias de processos qu\u00edmicos. Tradu\u00e7\u00e3o de Hor\u00e1cio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.Revistas:Brazilian Journal of Chemical Engineering, S\u00e3o Paulo, SP: Brazilian Society of Chemical Engineering, v. 11, n. 1, 1995-;Qu\u00edmica & Derivados, S\u00e3o Paulo, SP: QD, v. 1, n. 1, nov. 1965-; Dispon\u00edvel em: http://www.quimica.com.br/pquimica/category/revista/Petr\u00f3leo & Energia, S\u00e3o Paulo, SP, v. 1, n. 1, ; Dispon\u00edvel em: http://www.petroleoenergia.com.br/petroleo/category/revista-petroleo-e-energia/Revista FACTO, Publica\u00e7\u00e3o da Associa\u00e7\u00e3o Brasileira das Ind\u00fastrias de Qu\u00edmica Fina, Biotecnologia e suas Especialidades, Rio de Janeiro, RJ, v. 1, n. 1; Dispon\u00edvel em: http://www.abifina.org.br/facto/\"\n  ]\n];\n\nconst body = context.document.body;\nlet totalReplaced = 0;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWildcards: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + before.slice(0, 60));\n  }\n\n  // Replace from the last match to the first so that earlier ranges stay valid\n  // even if insertText were ever to shift sibling range anchors.\n  for (let i = results.items.length - 1; i >= 0; i--) {\n    results.items[i].insertText(after, \"Replace\");\n  }\n  totalReplaced += results.items.length;\n  await context.sync();\n}\n\nawait context.sync();\nreturn totalReplaced;\n", "ps1": "$d = $word.ActiveDocument\n\n# Word COM constants (wdFindContinue / wdReplaceAll)\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n  @('Ativa\u00e7\u00e3o: 01/01/2022', 'Ativa\u00e7\u00e3o: 01/01/2024'),\n  @('Curso (semestre ideal): EQD (9), EQN (10)', 'Curso (semestre ideal): EQN (10)'),\n  @('Conferir aos alunos uma no\u00e7\u00e3o pr\u00e1tica das aplica\u00e7\u00f5es, \u00e0 escala industrial, de processos e produtos estudados em disciplinas de qu\u00edmica org\u00e2nica e ainda uma vis\u00e3o global das mat\u00e9rias primas mais importantes na \u00e1rea da ind\u00fastria qu\u00edmica de base org\u00e2nica.', 'Conferir aos alunos aplica\u00e7\u00f5es em escala industrial, de processos e produtos estudados em disciplinas de qu\u00edmica org\u00e2nica e uma vis\u00e3o global das mat\u00e9rias primas mais importantes na \u00e1rea da ind\u00fastria qu\u00edmica de base org\u00e2nica.'),\n  @('Check the students a practical notion of applications on an industrial scale, processes and products studied in organic chemistry disciplines and also an overview of the most important raw materials in the chemical industry of organic base', 'Check the students students with applications on an industrial scale of processes and products studied in organic chemistry disciplines and an overview of the most important raw materials in the chemical industry of organic base.'),\n  @('Petr\u00f3leo, G\u00e1s Natural e Petroqu\u00edmica; Qu\u00edmica Fina; Processos Unit\u00e1rios Org\u00e2nicos: Nitra\u00e7\u00e3o; Esterifica\u00e7\u00e3o; Alquila\u00e7\u00e3o e Acila\u00e7\u00e3o; Hidrogena\u00e7\u00e3o; Sulfona\u00e7\u00e3o/Sulfata\u00e7\u00e3o; Oxida\u00e7\u00e3o.', 'Qu\u00edmica Fina,Petr\u00f3leo, G\u00e1s Natural e Petroqu\u00edmica; Processos Unit\u00e1rios Org\u00e2nicos:Alquila\u00e7\u00e3o e Acila\u00e7\u00e3o; Hidrogena\u00e7\u00e3o e Desidrogena\u00e7\u00e3o,Halogena\u00e7\u00e3o, Esterifica\u00e7\u00e3o; Sulfona\u00e7\u00e3o/Sulfata\u00e7\u00e3o; Oxida\u00e7\u00e3o.'),\n  @('Oil, Natural Gas and Petrochemicals; Fine Chemistry; Organic Unit Processes: Nitration; Esterification; Alkylation and Acylation; Hydrogenation; Sulphonation/Sulfation; Oxidation.', 'Fine Chemistry,Oil, Natural Gas and Petrochemicals; Organic Unit Processes: Alkylation and Acylation; Hydrogenation,Halogenation, Esterification; Sulphonation/Sulfation; Oxidation.'),\n  @('Petr\u00f3leo, G\u00e1s Natural e Petroqu\u00edmica\u037e 2- Qu\u00edmica Fina: Caracter\u00edsticas, Qu\u00edmica Fina X Qu\u00edmica de Base, Principais Segmentos (Defensivos Agr\u00edcolas, F\u00e1rmacos, Catalisadores, Corantes e Pigmentos, Especialidades); 3- Processos Unit\u00e1rios Org\u00e2nicos: 3.1- Nitra\u00e7\u00e3o; 3.2- Esterifica\u00e7\u00e3o; 3.3- Alquila\u00e7\u00e3o e Acila\u00e7\u00e3o; 3.4- Hidrogena\u00e7\u00e3o; 3.5- Sulfona\u00e7\u00e3o/Sulfata\u00e7\u00e3o; 3.6- Oxida\u00e7\u00e3o.', '1-Qu\u00edmica Fina: 1.1- Caracter\u00edsticas, 1.2- Principais Segmentos (Defensivos Agr\u00edcolas, F\u00e1rmacos, Catalisadores, Corantes e Pigmentos, Especialidades), 1.3- Qu\u00edmica Fina X Qu\u00edmica de Base, 2- Petr\u00f3leo, G\u00e1s Natural e Petroqu\u00edmica; 3- Processos Unit\u00e1rios Org\u00e2nicos: 3.1- Alquila\u00e7\u00e3o e Acila\u00e7\u00e3o; 3.2- Hidrogena\u00e7\u00e3o e Desidogena\u00e7\u00e3o;3.2.1 Processos Oxo, 3.2.2- Amino;3.3- Halogena\u00e7\u00e3o; 3.4- Esterifica\u00e7\u00e3o; 3.5- Sulfona\u00e7\u00e3o/Sulfata\u00e7\u00e3o; 3.6- Oxida\u00e7\u00e3o.'),\n  @('1- Oil, Natural Gas and Petrochemicals 2- Fine Chemistry: Characteristics, Fine Chemistry vs. Basic Chemistry, Main Segments (Pesticides, Drugs, Catalysts, Dyes and Pigments, Specialties); 3- Organic Unit Processes: 3.1- Nitration; 3.2- Esterification; 3.3- Alkylation and Acylation; 3.4- Hydrogenation; 3.5- Sulphonation/Sulfation; 3.6- Oxidation.', '1- Fine Chemicals: 1.1- Characteristics, 1.2- Main Segments (Agricultural Defensives, Pharmaceuticals, Catalysts, Dyes and Pigments, Specialties), 1.3- Fine Chemicals X Basic Chemicals, 2- Oil, Natural Gas and Petrochemicals; 3- Unitary Processes Organic: 3.1- Alkylation and Acylation; 3.2- Hydrogenation and Dehydrogenation; 3.2.1 Oxo processes, 3.2.2- Amino; 3.3- Halogenation; 3.4- Esterification; 3.5- Sulfonation/Sulfation; 3.6- Oxidation.'),\n  @('Provas em sala, entrega de exerc\u00edcios ou casos pr\u00e1ticos elaborados fora de sala de aula.', 'A nota (NOTA) ser\u00e1 composta por uma destas op\u00e7\u00f5es: prova em sala, apresenta\u00e7\u00f5es em sala, entrega de exerc\u00edcios ou casos pr\u00e1ticos elaborados fora de sala de aula. A estas op\u00e7\u00f5es ser\u00e1 incorporado,para cada aluno, seu respectivo percentual defrequ\u00eancia no c\u00e1lculo da nota final (NF), conforme a f\u00f3rmula explicitada abaixo:NF = NOTA x %FREQ.'),\n  @('Frequ\u00eancia m\u00ednima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recupera\u00e7\u00e3o.', 'Frequ\u00eancia m\u00ednima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita aplica\u00e7\u00e3o de prova escrita de recupera\u00e7\u00e3o valendo 10,00 pontos.'),\n  @('Ullmann\u2019s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim \u037e New York : WileyVCH, 2011.Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.Manual econ\u00f4mico da ind\u00fastria qu\u00edmica - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Cama\u00e7ari: CEPED, 2007.Shreve, R. Norris\u037e BRINK JR., J. A. Ind\u00fastrias de processos qu\u00edmicos. Tradu\u00e7\u00e3o de Hor\u00e1cio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.T.W. Graham Solomons, Craig B. Fryhle Hoboken, NJ. Organic chemistry; John Wiley, 9th ed; c2008.Revistas:Brazilian Journal of Chemical Engineering, S\u00e3o Paulo, SP: Brazilian Society of Chemical Engineering, v. 11, n. 1, 1995-;Qu\u00edmica & Derivados, S\u00e3o Paulo, SP: QD, v. 1, n. 1, nov. 1965-; Dispon\u00edvel em: http://www.quimica.com.br/pquimica/category/revista/BiodieselBR. Curitiba: BiodieselBR, v.1, n.1, out/dez.2007 -;Petr\u00f3leo & Energia, S\u00e3o Paulo, SP: , v. 1, n. 1, ; Dispon\u00edvel em: http://www.petroleoenergia.com.br/petroleo/category/revista-petroleo-e-energia/.', 'Livros:Ullmann\u2019s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim ; New York : WileyVCH, 2011.Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.Manual Econ\u00f4mico da Ind\u00fastria Qu\u00edmica - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Cama\u00e7ari: CEPED, 2007.Shreve, R. Norris; BRINK JR., J. A. Ind\u00fastrias de processos qu\u00edmicos. Tradu\u00e7\u00e3o de Hor\u00e1cio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.Revistas:Brazilian Journal of Chemical Engineering, S\u00e3o Paulo, SP: Brazilian Society of Chemical Engineering, v. 11, n. 1, 1995-;Qu\u00edmica & Derivados, S\u00e3o Paulo, SP: QD, v. 1, n. 1, nov. 1965-; Dispon\u00edvel em: http://www.quimica.com.br/pquimica/category/revista/Petr\u00f3leo & Energia, S\u00e3o Paulo, SP, v. 1, n. 1, ; Dispon\u00edvel em: http://www.petroleoenergia.com.br/petroleo/category/revista-petroleo-e-energia/Revista FACTO, Publica\u00e7\u00e3o da Associa\u00e7\u00e3o Brasileira das Ind\u00fastrias de Qu\u00edmica Fina, Biotecnologia e suas Especialidades, Rio de Janeiro, RJ, v. 1, n. 1; Dispon\u00edvel em: http://www.abifina.org.br/facto/'),\n)\n\nforeach ($pair in $replacements) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $found = $find.Execute($old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new, $wdReplaceAll)\n  if (-not $found) {\n    throw \"No match found for: \" + $old.Substring(0, [Math]::Min(60, $old.Length))\n  }\n}\n"}
